$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: values sourced from former row 14
$ws.Range("D2").Value = 44435
$ws.Range("J2").Value = 25
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 14000
$ws.Range("P2").Value = 467

# Row 3: values sourced from former row 15
$ws.Range("D3").Value = 44435

# Row 4: values sourced from former row 8
$ws.Range("D4").Value = 44446
$ws.Range("J4").Value = 25
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 14000
$ws.Range("P4").Value = 467

# Row 5: values sourced from former row 6
$ws.Range("D5").Value = 44467
$ws.Range("J5").Value = 35

# Row 6: values sourced from former row 3
$ws.Range("D6").Value = 44432
$ws.Range("J6").Value = 25
$ws.Range("K6").Value = 14000
$ws.Range("L6").Value = 14000
$ws.Range("M6").Value = 14000
$ws.Range("O6").Value = "Provincia del Elquí"
$ws.Range("P6").Value = 467

# Row 7: values sourced from former row 11
$ws.Range("D7").Value = 44841
$ws.Range("J7").Value = 45
$ws.Range("K7").Value = 12000
$ws.Range("L7").Value = 12000
$ws.Range("M7").Value = 12000
$ws.Range("P7").Value = 400

# Row 8: values sourced from former row 12
$ws.Range("D8").Value = 44841
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 45
$ws.Range("K8").Value = 10000
$ws.Range("L8").Value = 10000
$ws.Range("M8").Value = 10000
$ws.Range("N8").Value = "$/caja 40 unidades"
$ws.Range("P8").Value = 250
$ws.Range("Q8").Value = 40

# Row 9: values sourced from former row 5
$ws.Range("D9").Value = 44449
$ws.Range("K9").Value = 12000
$ws.Range("L9").Value = 12000
$ws.Range("M9").Value = 12000
$ws.Range("P9").Value = 400

# Row 10: values sourced from former row 7
$ws.Range("D10").Value = 44421
$ws.Range("J10").Value = 25
$ws.Range("L10").Value = 16000
$ws.Range("M10").Value = 15400
$ws.Range("P10").Value = 513

# Row 11: values sourced from former row 16
$ws.Range("D11").Value = 44425
$ws.Range("J11").Value = 35
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 14000
$ws.Range("M11").Value = 14000
$ws.Range("P11").Value = 467

# Row 12: values sourced from former row 9
$ws.Range("D12").Value = 44460
$ws.Range("I12").Value = "Primera"
$ws.Range("K12").Value = 13000
$ws.Range("L12").Value = 13000
$ws.Range("M12").Value = 13000
$ws.Range("N12").Value = "$/caja 30 unidades"
$ws.Range("P12").Value = 433
$ws.Range("Q12").Value = 30

# Row 14: values sourced from former row 4
$ws.Range("D14").Value = 44453
$ws.Range("J14").Value = 50
$ws.Range("K14").Value = 12000
$ws.Range("L14").Value = 12000
$ws.Range("M14").Value = 12000
$ws.Range("P14").Value = 400

# Row 15: values sourced from former row 10
$ws.Range("D15").Value = 44418
$ws.Range("J15").Value = 30
$ws.Range("K15").Value = 15000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 15000
$ws.Range("O15").Value = "Provincia de Limarí"
$ws.Range("P15").Value = 500

# Row 16: values sourced from former row 2
$ws.Range("D16").Value = 44474
$ws.Range("J16").Value = 45
$ws.Range("K16").Value = 10000
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = 10000
$ws.Range("P16").Value = 333
